$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds 4 years (2014-2017) of monthly data in rows 2..49 (12 rows/year).
# Each year's block needs reordering so Oct/Nov/Dec come first, followed by Jan..Sep
# (i.e. within each year, rotate the 12 months so the row order becomes
#  Oct, Nov, Dec, Jan, Feb, Mar, Apr, May, Jun, Jul, Aug, Sep).

$firstRow = 2
$lastRow = 49
$rowsPerYear = 12

# Read all existing rows (A:D) into memory first.
$data = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $c = $ws.Cells.Item($r, 3).Value()
    $d = $ws.Cells.Item($r, 4).Value()
    $data += ,@($a, $b, $c, $d)
}

$yearCount = $data.Length / $rowsPerYear

# Build the reordered list: for every year-block, move the last 3 rows
# (Oct, Nov, Dec) ahead of the first 9 rows (Jan..Sep).
$newData = @()
for ($yearIdx = 0; $yearIdx -lt $yearCount; $yearIdx++) {
    $base = $yearIdx * $rowsPerYear
    $block = $data[$base..($base + $rowsPerYear - 1)]
    $tail = $block[9..11]
    $head = $block[0..8]
    $newData += $tail
    $newData += $head
}

# Write the reordered values back into the same A:D range.
for ($i = 0; $i -lt $newData.Length; $i++) {
    $r = $firstRow + $i
    $row = $newData[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
}
